$wb = $excel.ActiveWorkbook

# Sheet "Basic Game rubric": mark row 2 ("Camera") as done (1) and add a note date
$ws1 = $wb.Worksheets.Item("Basic Game rubric")
$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = "klaar op 12/03/'22"

# Sheet "Game extras": mark row 4 ("Load level from a file") as done (1)
$ws2 = $wb.Worksheets.Item("Game extras")
$ws2.Range("B4").Value = 1

# Update the active selection on each sheet to match the saved view state
$ws2.Activate()
$ws2.Range("B13").Select()

$ws1.Activate()
$ws1.Range("C6").Select()
